# Update the "run the code / check accuracy" bullet on the
# "任务：实现推理" slide so that it calls out the 100% accuracy figure
# explicitly, splitting it into three runs:
#   1) "请每个同学都运行代码，查看推理正确率是否接近"
#   2) "100%"
#   3) "？"

$p = $ppt.ActivePresentation

$needle   = "正确率收敛了吗"
$runText1 = "请每个同学都运行代码，查看推理正确率是否接近"
$runText2 = "100%"
$runText3 = "？"

$targetPara  = $null
$targetRange = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi)
                if ($para.Text -like "*$needle*") {
                    $targetPara  = $para
                    $targetRange = $tr
                }
            }
        }
    }
}

if ($targetPara -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    $textRange = $targetRange
    $startPos  = $targetPara.Start
    $oldLen    = $targetPara.Length

    # Replace the whole paragraph with the final, combined text first.
    $whole = $textRange.Characters($startPos, $oldLen)
    $whole.Text = $runText1 + $runText2 + $runText3

    # Re-assert the text of each segment individually so the engine
    # materialises three distinct runs at the exact boundaries we need.
    $run1 = $textRange.Characters($startPos, $runText1.Length)
    $run1.Text = $runText1

    $run2 = $textRange.Characters($startPos + $runText1.Length, $runText2.Length)
    $run2.Text = $runText2

    $run3 = $textRange.Characters($startPos + $runText1.Length + $runText2.Length, $runText3.Length)
    $run3.Text = $runText3

    Write-Host "Updated paragraph successfully"
}
